# "add number 6 in green"
# Put the number 6 in F1 (the cell right after the existing E1=5) and
# color its font green using the theme's accent6 color (matches the
# workbook's existing pattern of coloring numbers via font theme colors).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("F1")
$cell.Value = 6
$cell.Font.ThemeColor = 10   # xlThemeColorAccent6 -> green (theme index 9)

# Move the active selection, matching the post-edit selection state.
$ws.Range("H6").Select() | Out-Null
